$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Property" sheet, column D = "Private". Set it to TRUE for the data rows
# (rows 2-6), matching the commit: "set property's private value as true".
$rng = $ws.Range("D2:D6")
$rng.Value = $true

# Mirror the author's selection after making the edit.
[void]$rng.Select()

# Extend the TRUE/FALSE list validation (already used on column F) to also
# cover the newly-boolean column D.
[void]$rng.Validation.Add(3, 1, 1, '"TRUE,FALSE"')
